$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $cols = @("B", "C", "D", "E", "F", "G")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-Rows 49 50
Swap-Rows 76 77
Swap-Rows 82 83
Swap-Rows 86 87
Swap-Rows 109 110
Swap-Rows 147 148
Swap-Rows 152 153
Swap-Rows 157 158
Swap-Rows 162 163
Swap-Rows 175 176
Swap-Rows 189 190
Swap-Rows 227 228
Swap-Rows 232 233
Swap-Rows 251 252
Swap-Rows 253 254
Swap-Rows 366 367
Swap-Rows 370 371
Swap-Rows 404 405
Swap-Rows 415 416
Swap-Rows 417 418
Swap-Rows 454 455
Swap-Rows 465 466
Swap-Rows 485 486
Swap-Rows 487 488
Swap-Rows 497 498
Swap-Rows 502 503
Swap-Rows 511 512
Swap-Rows 537 538
Swap-Rows 541 542
Swap-Rows 610 611
Swap-Rows 632 633
Swap-Rows 778 779
Swap-Rows 782 783
Swap-Rows 805 806
Swap-Rows 807 808
Swap-Rows 831 832
Swap-Rows 833 834
Swap-Rows 861 862
Swap-Rows 872 873
Swap-Rows 878 879
Swap-Rows 884 885
Swap-Rows 887 888
Swap-Rows 902 903
Swap-Rows 946 947
